$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Fill in the new time-log entry on row 11 ---
# Date column (A): 11/11/2014 (serial 41954)
$ws.Range("A11").Value = 41954

# Time In (B) / Time Out (C) -- Time Out is entered first so that the
# shared-string table indexes these in the same order as the target workbook
$ws.Range("C11").Value = "7:00pm"
$ws.Range("B11").Value = "8:00am"

# Time spent in minutes (D)
$ws.Range("D11").Value = 720

# General category (E)
$ws.Range("E11").Value = "Front End Programming"

# Description (F)
$ws.Range("F11").Value = "Created login, register buttons, modal and functionality"

# --- Update the view state: scroll down a couple rows and move the
#     active selection to A12 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("A12").Select() | Out-Null
